$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor (name unchanged, values updated)
$ws.Range("B3").Value = 0.9930953655259126
$ws.Range("C3").Value = 0.9916337801413614
$ws.Range("D3").Value = 0.9729820265271649

# Row 4 - renamed from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9854741089323991
$ws.Range("C4").Value = 0.9847295918877119
$ws.Range("D4").Value = 0.9196670521670267

# Row 5 - renamed from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9971513494778254
$ws.Range("C5").Value = 0.9948672209730874
$ws.Range("D5").Value = 0.9910473861522036
